# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-11 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new K value (column G)
$kValues = @{
    2  = 0
    3  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
